$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(29, 1).Value = 45873
$ws.Cells.Item(29, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(29, 2).Value = 6329.93994140625
$ws.Cells.Item(29, 3).Value = 6330.68994140625
$ws.Cells.Item(29, 4).Value = 6271.7099609375
$ws.Cells.Item(29, 5).Value = 6271.7099609375
$ws.Cells.Item(29, 6).Value = 4842580000
$ws.Cells.Item(29, 7).Value = 0.014737100330916
$ws.Cells.Item(29, 8).Value = 1
$text29 = 'nintendo break sale record did delayed gta release save switch fedex corporation fdxs ceo is a hitter say jim cramer inside gamesxs ai revolution can india become the next global gamechanger intel corporation intc taught me not to trust ai chatbots say jim cramer ipo stock of the week china stock leader atour test buy point amid market fall the procter gamble company pg it just not great say jim cramer espn obtains redzone other medium asset from nfl in blockbuster deal per report this week in ai key to agentic ai adoption plus ai action plan from u china mcdonalds plan major ai expansion to streamline customer order executive say paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame oil major beat back billion fine in kazakhstan paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame the outdated management mindset holding your business back paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame i worked in the building where people were killed here how to make your building more secure paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame columbia is selling a stylish quilted jacket for and reviewer say it is wellmade paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame top action comedy anime series in paula deen ha abruptly closed the savannah restaurant that launched her to food network fame paula deen ha abruptly closed the savannah restaurant that launched her to food network fame already hit your k limit in here what to do next cadence design system inc cdnss number show nvidia is in great shape in china say jim cramer alnylams amvuttra make powerful entrance eye firstline dominance president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report nvidia corporation nvda could be in great shape in china say jim cramer the boeing company ba a horse it ceo ha done a great job say jim cramer doj identifies tariff fraud a enforcement priority waste management inc wms ceo is doing a remarkable job say jim cramer merck co inc mrk need to move on from talking about it animal division say jim cramer starbucks corporation sbux dont worry about luckin coffee say jim cramer sarah jessica parker announces emotional and poetic end to and just like that carrier global corporation carr delivered a good quarter say jim cramer iphone pro launch in september new colour a pro chip expected india price and more otis worldwide corporation otis had a tough quarter china wa horrible say jim cramer a apple watch series is just one click away at stacksocial celestica inc cl is up huge say jim cramer ai evolution prompt sense of urgency among cfo oracle jury say tesla wa partly to blame for fatal crash trump issue nuclear threat to russiaand wall street fear gauge go vertical crypto founder take note these sec roundtable could decide the rule that govern you sugar price pressured a brazil ramp up sugar production this week in stablecoins building nextgen rail for enterprise finance tesla must pay million in damage after fatal autopilot crash jury say cnbc miami jury find tesla partly liable in autopilot crash award million in punitive damage nbc news bond etf could catch a tailwind from job data tesla tsla is found liable in fatal autopilot crash ha to pay million the latest u stock market and global trade partner react to trump new tariff no one can be that wrong trump fire labor statistic bos after weak revised job data the green revolution anand chandra cofounder executive director aryaag on fixing india agri gap with tech cocoa price sharply lower on the outlook for adequate supply what ghislaine maxwell move to a lowersecurity federal prison camp could signal reddit want to replace google a the internet top search engine investor are all in ct bond commission approves m to lower electric bill republican say it not enough coffee price plunge a supply concern ease the latest business and u trading partner worldwide respond to trump new tariff the latest business and u trading partner worldwide respond to trump new tariff nextnav and oscilloquartz partner to advance g based timing solution for critical infrastructure resmed analyst increase their forecast after strong q earnings hedge fund and insider trading news william heard warren buffett ray dalio bill ackman rokos capital management alden global capital clear channel outdoor holding inc cco and more here are reason we are not going to buy stock in friday nosediving market this is why coinbase stock coin is cratering today tipranks phoenix az lemonade stand squeeze success expands into western state with over million in sale top cryptos to join now for that could deliver big return a new bowling center wa big news in west valley bowl face endofmonth deadline corporation for public broadcasting will shutter after trump funding cut the best trader in congress how to track politician stock trade a dallas bank is confident it loss will pay off best gold stock to watch in midtiers and junior beat the major online creator led by mrbeast and mark rober want to raise million for clean water access online creator led by mrbeast and mark rober want to raise million for clean water access online creator led by mrbeast and mark rober want to raise million for clean water access former chairman of savannah economic development dy jury order tesla to pay m in autopilot crash case opening it up to other costly lawsuit jury order tesla to pay m in autopilot crash case opening it up to other costly lawsuit online creator led by mrbeast and mark rober want to raise million for clean water access tesla found partially liable for fatal crash involving driverassistance technology jury order million in damage the washington post india secret weapon for defence sovereignty through ai and zerochina tech trump injects a new dose of uncertainty in tariff a he push start date back to aug trump injects a new dose of uncertainty in tariff a he push start date back to aug these analyst boost their forecast on microstrategy following q result trump injects a new dose of uncertainty in tariff a he push start date back to aug unit moving and portable storage return a primary sponsor for timmy hill at watkins glen nascar truck series race tech stock lead selloff a soft labor data fuel fed cut bet the eye of the hurricane why the u job market ha soured economist say here what could get more expensive from trump massive tariff hike here what could get more expensive from trump massive tariff hike trump seek to leave his goldplated stain on the white house beatrice police department hosting national night out nicole leone former executive director erie humane society ha september trial date president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report president trump seek firing of government official overseeing job data after weak employment report north carolina u senate race turn into a battle for the middle class sarah jessica parker announces emotional and poetic end to and just like that ghislaine maxwell join famous prisoner in move to texas prison camp wave of fake credential spark political fallout in spain gop congressman suggests hed lose money and go broke if adviser stop stock trade gop congressman suggests hed lose money and go broke if adviser stop stock trade gop congressman suggests hed lose money and go broke if adviser stop stock trade jeffrey epstein former girlfriend ghislaine maxwell is transferred to a prison camp in texas jeffrey epstein former girlfriend ghislaine maxwell is transferred to a prison camp in texas jeffrey epstein former girlfriend ghislaine maxwell is transferred to a prison camp in texas jeffrey epstein former girlfriend ghislaine maxwell is transferred to a prison camp in texas gop congressman suggests hed lose money and go broke if adviser stop stock trade gop congressman suggests hed lose money and go broke if adviser stop stock trade gop congressman suggests hed lose money and go broke if adviser stop stock trade corporation for public broadcasting which fund pb npr to close after federal aid cut nyc subway nightmare flood trap commuter in filth and failure the latest u stock market and global trade partner react to trump new tariff what ghislaine maxwell move to a lowersecurity federal prison camp could signal what new public comment rule mean for speaking at local meeting in bismarckmandan area corporation for public broadcasting will shutter after trump funding cut congressman view here the truth about impact of beautiful bill on st louis county trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev state college rally planned a another local health system cut youth trans care more american couple are turning to italy dolce vita in a quest for memorable wedding city say thursday odor came from construction of city new industrial water plant trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia cult of cruelty ghastly photo video of emaciated israeli hostage released by gaza terrorist here what could get more expensive from trump massive tariff hike red state move to slash medicaid funding for planned parenthood trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia trump announces nuclear submarine positioned near russia trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev trump order u nuclear sub repositioned over statement from exrussian leader medvedev guest columnist we must empower support girl and boy trump order nuclear submarine repositioned after provocative statement from russia editorial maxwell talk giving more life to a conspiracy shribman who a hero in the age of trump brown fulfilling the promise of the ada a trophy hunter killed a lion in zimbabwe that wa part of a research project sparking anger pensacola man traveled to orlando stealing smokey bear sign along the way official say kamala harris insufferably cackle over the government spying on american through wireless earpods western country speak of a future palestinian state a the nightmare unfolding in gaza worsens western country speak of a future palestinian state a the nightmare unfolding in gaza worsens western country speak of a future palestinian state a the nightmare unfolding in gaza worsens western country speak of a future palestinian state a the nightmare unfolding in gaza worsens western country speak of a future palestinian state a the nightmare unfolding in gaza worsens trump say he deploying two nuclear sub after provocative comment from exrussian leader democrat are struggling to block trump texas power grab the healing power of sunlight office cleaner who came facetoface with the manhattan gunman detail his deadly rampage corporation for public broadcasting funder of npr and pb say it will end operation within month after federal budget cut dozen of protester swarm chuck schumer kirsten gillibrands nyc office to let gaza live fourth circuit rule unanimously in favor of employee in two case involving employer allegation of breach of restrictive covenant city county official talk about animal control a texas researcher wa held at an airport for over a week now he face deportation a texas researcher wa held at an airport for over a week now he face deportation a texas researcher wa held at an airport for over a week now he face deportation a texas researcher wa held at an airport for over a week now he face deportation a texas researcher wa held at an airport for over a week now he face deportation trump plan to build million white house party room america rising part a federal judge just called out the trump administration for lying to the supreme court k nieth to retire fbi official reportedly redacted reference to trump from epstein file trump positioning two nuclear submarine near russia following russian official threat of war how immigration legislation is dangerous to immigrant judge pause removal of k migrant from honduras nicaragua and nepal in overthetop ruling democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan democrat launching summer blitz to press republican on trump spending plan trump weighs in on american eagle sydney sweeney ad shiba inu and cardano set to breakout a one new token set sight on gain top cryptos to buy in trimiq weight loss official launched in uk and ireland with green tea probiotic complex free spin no deposit bonus casino online real money casino no deposit by wild casino the shock job report set off this recession alert and hold fresh clue that ai may be boosting unemployment jpmorgan say fortune crypto presale that could make you a millionaire in dogecoin ada xrp cant match ozak ai potential arialief for neuropathy supplement official launch of a new hope for sciatic nerve discomfort relief the ig choice ha gone off the rail reader commentary what to know if youre at risk of having your wage garnished over student loan debt m raised and climbing ruvi ai ruvi called ripple xrp a coinmarketcap listing made it token one of the highest in demand why parent are suing snapchat over fentanyl death tea encouraged it user to spill then the apps data got leaked npr federal reserve governor resigns creating vacancy for trump federal reserve governor resigns creating vacancy for trump upicom bragar eagel squire pc is reminds investor of class action against neogen alto bragar eagel squire pc is reminds investor of class action against neogen alto rxsight and replimune and encourages investor to contact the firm audizen drop officially launched natural tinnitus support inspired by audifort europe ev sale report tesla ha plus win in a ev share market matsato knife official launch in australia usa redefining precision and power in the kitchen warren buffett dividend stock you can buy and hold forever crushed q earnings luka doncic agrees to a threeyear million contract extension through with the la lakers the financial swamp of banking regulator who decide which american get their bank account shut down sam altman compare openai to the manhattan projectand he not joking about the risk the democrat have a mamdani problem editorial toncoin at ondo cross while cold wallet entry crypto reward outclass both weekend gaming best ipad game to try in little pepe lilpepe price prediction for show strong potential to beat solana sol and cardano ada sluggish summer in sin city la vega see big drop in tourist ap business summarybrief at pm edt ap business summarybrief at pm edt from lao to brazil trump tariff leave a lot of loser but even the winner will pay a price from lao to brazil trump tariff leave a lot of loser but even the winner will pay a price from lao to brazil trump tariff leave a lot of loser but even the winner will pay a price from lao to brazil trump tariff leave a lot of loser but even the winner will pay a price texas panel advance redrawn congressional map that could take democratic seat your money exit a business or equity position without regret pepes chart turn positive litecoin gather whale support blockdag expands web reach with devs a day workweek isnt good just for team moral it also good for the bottom line teleflex raise outlook following strong q and biotronik acquisition tucson trolley tour operator seek to expand reach former trump prosecutor jack smith probed for alleged election interference u hiring slowed sharply over the summer the wall street journal mstr news michael saylor explains why strategy strc preferred stock is the firm iphone moment coindesk harvard say it will comply with dhs request to hand over employment form missed binance coin bnb expert say ruvi ai ruvi audited token is next millionaire maker a coinmarketcap listing spark high demand what to know if youre at risk of having your wage garnished over student loan debt how to safely download and install game on your laptop shib find support ondo hold while cold wallet price gas rebate make it a top crypto to watch little pepe lilpepe gather momentum with gain ahead while cardano ada nears run the trump administration take a very orwellian turn how this week avalanche of news from washington to wall street kept investor guessing westchester capital management inc trim stake in valero energy corporation nysevlo united super pty ltd in it capacity a trustee for the construction building union superannuation fund make new investment in valero energy corporation nysevlo decadeold bitcoin wallet emptied million in btc moved amid market correction how florida quietly surpassed california in solar growth ukrainian longrange drone hit multiple russian refinery factory overnight journointerior designer jim acosta guarantee trump wh ballroom is going to look like st americu appoints new community relation manager mitsubishi ufj trust banking corp purchase share of ge aerospace nysege j safra sarasin holding ag purchase share of thermo fisher scientific inc nysetmo mitsubishi ufj trust banking corp sell share of thermo fisher scientific inc nysetmo banco santander sa ha million stake in thermo fisher scientific inc nysetmo the walt disney company nysedis share bought by bank of nova scotia trust co the walt disney company nysedis stake lifted by r squared ltd the walt disney company nysedis share acquired by j safra sarasin holding ag the walt disney company nysedis share bought by j safra sarasin holding ag banco santander sa increase stock holding in the walt disney company nysedis this ai is set to flip every major altcoinare you in figma ipo surprise winner is a charity with million sharesand a famous backstory that sparked a bitter feud over an oil fortune decade ago fortune the hidden signal in oil market moor cabot inc ha million stock holding in the goldman sachs group inc nysegs moor cabot inc buy share of the goldman sachs group inc nysegs hartford investment management co cut stock position in zoetis inc nysezts how florida quietly surpassed california in solar growth cnbc lesson from a decade of empowering woman entrepreneur to lead boldly best trucking bookkeeping service three stock to consider today blockdags no vesting pas with roi end soon eth eye jup drop just of smallbusiness owner say the child care crisis is a critical problem here why pepeto price forecast will this presale meme coin surpass dogecoin shiba inu and ethereum editorial beyond the number study into subminimum must consider all ramification thing to know about palantir pltr before it report q earnings best trucking bookkeeping service it super disrespectful wnba player exasperated after sex toy thrown onto court for second time in a week claude cutoff anthropic revoke openais api access ahead of gpt launch malaga cove capital llc sell share of abbott laboratory nyseabt bayforest capital ltd ha stake in fidelity national financial inc nysefnf vestcor inc buy share of fidelity national financial inc nysefnf modera wealth management llc sell share of arista network inc nyseanet bank of nova scotia sell share of procter gamble company the nysepg bank of nova scotia increase stock holding in arista network inc nyseanet aviso wealth management acquires share of procter gamble company the nysepg aviso wealth management cut position in the boeing company nyseba j safra sarasin holding ag purchase share of abbott laboratory nyseabt pathway financial adviser llc purchase share of eli lilly and company nyselly aviso wealth management lower stock position in the boeing company nyseba procter gamble company the nysepg share bought by j safra sarasin holding ag what are the best dividend growth stock i should invest in right now hims alert bragar eagel squire pc announces that a class action lawsuit ha been steve witkoff meet with hostage family in tel aviv susan estrich the harvard fight money v principle jet holiday ad audio used in ice deportation video shes so cringe christine pelosi tease that shes running but not the way youd think video from the dream of dr king to the nightmare that is donald trump opinion from the dream of dr king to the nightmare that is donald trump opinion from the dream of dr king to the nightmare that is donald trump opinion guam d printing facility to enhance navy maintenance capability the ig choice ha gone off the rail reader commentary trump motif can be opaque i asked a feminist economist to bare them india stand firm on russian oil import despite trump sanction threat pedaling dei cycling weekly demand bike community include fat black woman or theyre racist the u military just launched it largest air exercise in year pennsylvania lawmaker say child care is broken but dont agree on how to fix it yet twice told tale city council attends team building workshop if trump save maxwell he will destroy everything he promised pope thrill hundred of thousand of young catholic at holy year youth festival pope thrill hundred of thousand of young catholic at holy year youth festival are republican doing enough to confirm trump nominee conservative critic say no texas panel advance redrawn congressional map that could take democratic seat trump wa asked about his press secretary performance his answer took a weird turn a i see it corvallis ha done nothing to help gazans how top democrat are already gearing up for online texas is moving forward with it radical redistricting plan america trucker need a place to park george mason university president keep job get a raise despite trump admins probe in which way is the bls biased donna brazile say nprpbsdefunding maga cannot compete with idea cue a hard trip over hillary right degree wrong time new graduate face discouraging odds of landing a job trump administration move to ban abortion care for veteran unspeakably cruel manhunt intensifies for suspect in montana bar mass shooting texas panel advance redrawn congressional map that could take democratic seat bill maher know trump losing credibility when these supporter smell bullst kevin oleary slam trump for whacking bls chief after disappointing job report you dont shoot the messenger trump administration halted civil right lawsuit targeting abuse of prisoner and the mentally ill it wa never about safety one week in here whats blocked by the uk orwellian online bill these three country dominate the fighter jet export marketfor now tiktok can shape america next generation and beijing know it former trump prosecutor jack smith probed for alleged election interference virginia end fiscal year with million surplus for next budget bonus editorial cartoon for aug other voice conserving enhancing agricultural resource with solar energy trump alum share chilling google message from before secondterm return lawfare at it finest texas house panel advance redrawn congressional map that would add more gop seat letter nd representative helping put our nation closer to fiscal crisis tracey weiss our ocean backyard coming soon opening the gate to cotonicoast dairy senate make progress in averting government shutdown much earlier than usual no piecemeal deal witkoff tell hostage family trump want full gaza agreement will the wnba approve record m deal to bring connecticut sun to boston proud democrat call into chris cuomos show admits she wa wrong about trump shes not the only one texas house committee advance gopfriendly map texas house committee advance gopfriendly map texas house committee advance gopfriendly map release the epstein file and end the coverup texas house committee advance gopfriendly map texas house committee advance gopfriendly map texas house committee advance gopfriendly map the trump administration take a very orwellian turn new photo reveal iconic white house rose garden paved over after trump makeover the nfls new kickoff rule is here to stay with a slight tweak what else is changing oversea chinese banking corp ltd raise position in ecolab inc nyseecl westchester capital management inc ha million stake in becton dickinson and company nysebdx ecolab inc nyseecl share purchased by mirabaud cie sa limited option for democrat to retaliate if texas republican redraw congressional map trump hit brazilian product with tariff over bolsonaro sign that someone you know is secretly broke do they apply to the people around you party poopers less than of american friendship cross political line texas house committee advance redrawn congressional map congressional redistricting map pass house committee pushing dems towards quorum break how trump became the new master of the senate luka doncic lakers agree to year million extension banco santander sa grows stock position in the walt disney company nysedis union pacific corporation nyseunp share sold by mitsubishi ufj trust banking corp moor cabot inc ha million stock holding in zoetis inc nysezts black swimmer teach others to help end history of aquatic segregation laduke let the sun shine so say former intelligence agent rolling stone claim the kremlin made up russiagate evidence can donald trump f diplomacy keep everyone happy laduke let the sun shine laduke let the sun shine in japan launched a second strike on pearl harborsort of laduke let the sun shine your view people have the right to sign recall petition laduke let the sun shine letter to the editor the generosity of kentuckian laduke let the sun shine laduke let the sun shine it super disrespectful wnba player exasperated after sex toy thrown onto court for second time in a week it trump economy now the latest financial number offer some warning sign aviso wealth management sell share of eli lilly and company nyselly j safra sarasin holding ag buy share of procter gamble company the nysepg steve witkoff meet with hostage family in tel aviv u court upholds order blocking indiscriminate targeting by immigration patrol u court upholds order blocking indiscriminate targeting by immigration patrol oversea chinese banking corp ltd reduces stock position in welltower inc nysewell here what we know about the anaconda shooter the victim and ongoing search effort banque cantonale vaudoise purchase new position in entergy corporation nyseetr the sherwinwilliams company nyseshw position lifted by banco santander sa lsv asset management sell share of phillips nysepsx profit drop at warren buffetts berkshire hathaway a it writes down it kraft heinz investment moor cabot inc sell share of thermo fisher scientific inc nysetmo trump tariff are making money that may make them hard to quit the new york time revolution camp unhinged white liberal spend weekend smashing junk in antitrump fury sen josh hawleys longshot stock trading ban bill advance rankles trump larimer county restaurant inspection made july st louis boeing worker again reject contract strike to begin at midnight former associate executive director file lawsuit against ahsaa boeings st louis union worker will strike on monday boeing defense union poised for it first strike since bloombergcom mahjong the latest networking opp for woman entrepreneur why analyst say this audited ai token is the next big thing coinmarketcap listing spark daily sale in the million luxury label lose their luster amid changing vibe gucci sale plunge and lvmh disappoints best crypto to buy now and secure early gain north texas father of win million at choctaw casino this cookware maker is bracing for steel tariff behind a wall of pot and pan nbc news opec country to boost oil production by barrel per day opec country to boost oil production by barrel per day opec country to boost oil production by barrel per day opec country to boost oil production by barrel per day opec country to boost oil production by barrel per day trump live update white house defends firing bls commissioner over weak job number the new york time opec country to boost oil production by barrel per day opec country to boost oil production by barrel per day white house official defend firing of labor official a critic warn of trust erosion clickbait and clean water youtubers are the new philanthropist palantir advanced micro walt disney pfizer zoominfo and more stock to watch this week barrons sen blackburn to introduce bill to root out embedded foreign interest free iq test online with instant result reliable accurate free realiqonline roll out new iq testing experience himerex pro this himerex pro uk app set new standard in aidriven trading with unmatched security and user approval ev lead vehicle category in california berkshire earnings key takeaway strong profit no buyback high cash level barrons steam and itchio pull hundred of adult game over alleged visa mastercard pressure pauline ferrandprevot win tour de france femmes becoming first french rider to win event since iam district member in st louis reject latest boeing offer strike for fair contract iam union bmw sue teen entrepreneur twice over alibaba trademark infringement top trader see ripple xrp hitting binance coin bnb and little pepe lilpepe when to buy for maximum roi starbucks to close pickup store in back to starbucks pivot citigroup launch stratum elite card to rival chase sapphire reserve disney world ferry crash into dock in stormy weather injuring rider u manufacturing resurgence foreign firm relocate amid tariff the fantastic four first step hold it lead atop the box office the fantastic four first step hold it lead atop the box office the fantastic four first step hold it lead atop the box office la vega unemployment at amid tourism drop and visitor de'
$ws.Cells.Item(29, 9).Value = $text29
$ws.Cells.Item(29, 10).Value = $text29
$ws.Cells.Item(29, 11).Value = 0.124
$ws.Cells.Item(29, 12).Value = 0.097
$ws.Cells.Item(29, 13).Value = 0.779
$ws.Cells.Item(29, 14).Value = 0.9995000000000001
$ws.Cells.Item(29, 15).Value = 189
$ws.Cells.Item(29, 16).Value = 151
$ws.Cells.Item(29, 17).Value = 0
$ws.Cells.Item(29, 18).Value = 348
$ws.Cells.Item(29, 19).Value = 224
$ws.Cells.Item(29, 20).Value = 473
$ws.Cells.Item(29, 21).Value = 258
$ws.Cells.Item(29, 22).Value = 112
$ws.Cells.Item(29, 23).Value = 71
$ws.Cells.Item(29, 24).Value = 177
